$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.577.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.489.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.31%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '657.39'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.45%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.418'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.486.58'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.34'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +11.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.204'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '97.475.70'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.17'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.137.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000256'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.76'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.488.25'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.49'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +11.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.00'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +15.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.503'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '526.03'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.34'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000198'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '96.67'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.48'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.670.17'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.33'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +15.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.83'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +16.47%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.188'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.589'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.21'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.82%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.85'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.75%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.25%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '514.47'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.911'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +11.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.36'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.72'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +6.32%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.94%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.33'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.61'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.23'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +13.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.48'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.57%  '
